$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 218.125
$ws.Range("I2").Value = 273.75
$ws.Range("K2").Value = 273.75
$ws.Range("M2").Value = -160.75
$ws.Range("H62").Value = 1237.2
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
$ws.Range("H65").Value = 1237.2
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
$ws.Range("H76").Value = 4861.5386
$ws.Range("I76").Value = 3042.8572
$ws.Range("K76").Value = 3042.8572
$ws.Range("M76").Value = -2727.8572
$ws.Range("H79").Value = 4861.5386
$ws.Range("I79").Value = 3042.8572
$ws.Range("K79").Value = 3042.8572
$ws.Range("M79").Value = -1950.8572
$ws.Range("H100").Value = 1539.7693
$ws.Range("I100").Value = 761
$ws.Range("J100").Value = 2026.5
$ws.Range("K100").Value = 761
$ws.Range("L100").Value = 2026.5
$ws.Range("M100").Value = -220
$ws.Range("N100").Value = -3108.5
$ws.Range("H106").Value = 102700
$ws.Range("I106").Value = 127250
$ws.Range("K106").Value = 127250
$ws.Range("M106").Value = -126619
$ws.Range("H116").Value = 2945.7932
$ws.Range("I116").Value = 2158.1904
$ws.Range("K116").Value = 2158.1904
$ws.Range("M116").Value = 1283.8096
$ws.Range("H129").Value = 1206.9286
$ws.Range("J129").Value = 2472.2727
$ws.Range("L129").Value = 7416.8181
$ws.Range("N129").Value = -17416.8181
$ws.Range("H135").Value = 853.1053000000001
$ws.Range("I135").Value = 734.1111
$ws.Range("K135").Value = 6606.9999
$ws.Range("M135").Value = -4071.9999
$ws.Range("H137").Value = 6061374.5
$ws.Range("I137").Value = 655.2
$ws.Range("J137").Value = 15385558
$ws.Range("K137").Value = 1965.6
$ws.Range("L137").Value = 46156674
$ws.Range("M137").Value = 584.3999999999999
$ws.Range("N137").Value = -46161774
$ws.Range("H138").Value = 5377861
$ws.Range("I138").Value = 6945687.5
$ws.Range("J138").Value = 2456.7856
$ws.Range("K138").Value = 20837062.5
$ws.Range("L138").Value = 7370.3568
$ws.Range("M138").Value = -20831922.5
$ws.Range("N138").Value = -17650.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1030.9412
$ws.Range("I2").Value = 1030
$ws.Range("J2").Value = 1035.3334
$ws.Range("K2").Value = 1030
$ws.Range("L2").Value = 1035.3334
$ws.Range("M2").Value = -917
$ws.Range("N2").Value = -1261.3334
$ws.Range("H32").Value = 9121.24
$ws.Range("I32").Value = 9219.031000000001
$ws.Range("J32").Value = 8704
$ws.Range("K32").Value = 9219.031000000001
$ws.Range("L32").Value = 8704
$ws.Range("M32").Value = -8932.031000000001
$ws.Range("N32").Value = -9278
$ws.Range("H61").Value = 9805188
$ws.Range("I61").Value = 10870898
$ws.Range("J61").Value = 657.2
$ws.Range("K61").Value = 10870898
$ws.Range("L61").Value = 657.2
$ws.Range("M61").Value = -10870686
$ws.Range("N61").Value = -1081.2
$ws.Range("H74").Value = 7576954
$ws.Range("I74").Value = 10205300
$ws.Range("J74").Value = 1134.3529
$ws.Range("K74").Value = 10205300
$ws.Range("L74").Value = 1134.3529
$ws.Range("M74").Value = -10204426
$ws.Range("N74").Value = -2882.3529
$ws.Range("H77").Value = 7576954
$ws.Range("I77").Value = 10205300
$ws.Range("J77").Value = 1134.3529
$ws.Range("K77").Value = 51026500
$ws.Range("L77").Value = 5671.7645
$ws.Range("M77").Value = -51022132
$ws.Range("N77").Value = -14407.7645
$ws.Range("H116").Value = 1030.9412
$ws.Range("I116").Value = 1030
$ws.Range("J116").Value = 1035.3334
$ws.Range("K116").Value = 1030
$ws.Range("L116").Value = 1035.3334
$ws.Range("M116").Value = 1264
$ws.Range("N116").Value = -5623.3334
$ws.Range("H132").Value = 7578223
$ws.Range("I132").Value = 10418963
$ws.Range("J132").Value = 2917.111
$ws.Range("K132").Value = 31256889
$ws.Range("L132").Value = 8751.332999999999
$ws.Range("M132").Value = -31254359
$ws.Range("N132").Value = -13811.333
$ws.Range("H136").Value = 9805188
$ws.Range("I136").Value = 10870898
$ws.Range("J136").Value = 657.2
$ws.Range("K136").Value = 32612694
$ws.Range("L136").Value = 1971.6
$ws.Range("M136").Value = -32610144
$ws.Range("N136").Value = -7071.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1030.9412
$ws.Range("I3").Value = 1030
$ws.Range("J3").Value = 1035.3334
$ws.Range("K3").Value = 1030
$ws.Range("L3").Value = 1035.3334
$ws.Range("M3").Value = -916
$ws.Range("N3").Value = -1263.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6412385
$ws.Range("I31").Value = 1840.8043
$ws.Range("J31").Value = 55559892
$ws.Range("K31").Value = 1840.8043
$ws.Range("L31").Value = 55559892
$ws.Range("M31").Value = -1545.8043
$ws.Range("N31").Value = -55560482
$ws.Range("H34").Value = 6412385
$ws.Range("I34").Value = 1840.8043
$ws.Range("J34").Value = 55559892
$ws.Range("K34").Value = 1840.8043
$ws.Range("L34").Value = 55559892
$ws.Range("M34").Value = -1638.8043
$ws.Range("N34").Value = -55560296
$ws.Range("H53").Value = 24366.666
$ws.Range("J53").Value = 24366.666
$ws.Range("L53").Value = 24366.666
$ws.Range("N53").Value = -25580.666
$ws.Range("H58").Value = 1040.1702
$ws.Range("I58").Value = 485.07318
$ws.Range("J58").Value = 4833.3335
$ws.Range("K58").Value = 485.07318
$ws.Range("L58").Value = 4833.3335
$ws.Range("M58").Value = -282.07318
$ws.Range("N58").Value = -5239.3335
$ws.Range("H132").Value = 11629652
$ws.Range("I132").Value = 14707323
$ws.Range("J132").Value = 2894.889
$ws.Range("K132").Value = 44121969
$ws.Range("L132").Value = 8684.667000000001
$ws.Range("M132").Value = -44119439
$ws.Range("N132").Value = -13744.667
$ws.Range("H136").Value = 1040.1702
$ws.Range("I136").Value = 485.07318
$ws.Range("J136").Value = 4833.3335
$ws.Range("K136").Value = 1455.21954
$ws.Range("L136").Value = 14500.0005
$ws.Range("M136").Value = 1094.78046
$ws.Range("N136").Value = -19600.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H132").Value = 6334554.5
$ws.Range("I132").Value = 2929.476
$ws.Range("J132").Value = 31265328
$ws.Range("K132").Value = 8788.428
$ws.Range("L132").Value = 93795984
$ws.Range("M132").Value = -6258.428
$ws.Range("N132").Value = -93801044
$ws.Range("H136").Value = 19236544
$ws.Range("I136").Value = 26317742
$ws.Range("J136").Value = 16144.286
$ws.Range("K136").Value = 78953226
$ws.Range("L136").Value = 48432.858
$ws.Range("M136").Value = -78950676
$ws.Range("N136").Value = -53532.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 32000
$ws.Range("I14").Value = 58000
$ws.Range("J14").Value = 19000
$ws.Range("K14").Value = 58000
$ws.Range("L14").Value = 19000
$ws.Range("M14").Value = -57832
$ws.Range("N14").Value = -19336
$ws.Range("H33").Value = 16213.444
$ws.Range("J33").Value = 16213.444
$ws.Range("L33").Value = 16213.444
$ws.Range("N33").Value = -16713.444
$ws.Range("H36").Value = 16213.444
$ws.Range("J36").Value = 16213.444
$ws.Range("L36").Value = 16213.444
$ws.Range("N36").Value = -16713.444
$ws.Range("H51").Value = 13800.25
$ws.Range("I51").Value = 7070
$ws.Range("J51").Value = 16043.667
$ws.Range("K51").Value = 7070
$ws.Range("L51").Value = 16043.667
$ws.Range("M51").Value = -6560
$ws.Range("N51").Value = -17063.667
$ws.Range("H81").Value = 1013
$ws.Range("I81").Value = 750
$ws.Range("J81").Value = 1108.6364
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 2217.2728
$ws.Range("M81").Value = -439
$ws.Range("N81").Value = -4339.272800000001
$ws.Range("H84").Value = 1013
$ws.Range("I84").Value = 750
$ws.Range("J84").Value = 1108.6364
$ws.Range("K84").Value = 7500
$ws.Range("L84").Value = 11086.364
$ws.Range("M84").Value = -2196
$ws.Range("N84").Value = -21694.364
$ws.Range("H132").Value = 1214.8857
$ws.Range("I132").Value = 1054.7593
$ws.Range("J132").Value = 1755.3125
$ws.Range("K132").Value = 3164.2779
$ws.Range("L132").Value = 5265.9375
$ws.Range("M132").Value = -634.2779
$ws.Range("N132").Value = -10325.9375
